# Documentos->Convenios: Se da nueva estructura al decreto del director o
# subrogante en cada uno de los borradores.
#
# Before:  ...consta en el ${directorDecreto}, que aprobó el Reglamento
#          Orgánico de los Servicios de Salud. La representación...
# After:   ...consta en el ${art8}Decreto N°140/04, del Ministerio de Salud
#          que aprobó el Reglamento Orgánico de los Servicios de Salud,
#          ${directorDecreto}. La representación...

$d = $word.ActiveDocument

# --- Step 1: drop the old trailing sentence run that followed the
# ${directorDecreto} placeholder ("... que aprobó el Reglamento Orgánico
# de los Servicios de Salud"), it is being folded into the new prefix
# text below.
$rngOld = $d.Content
$foundOld = $rngOld.Find.Execute(', que aprobó el Reglamento Orgánico de los Servicios de Salud', `
    $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($foundOld) {
    $rngOld.Delete()
}

# --- Step 2: expand the ${directorDecreto} placeholder run so it carries
# the new lead-in text plus itself at the end. Setting .Text on the found
# range keeps the run's existing formatting (Calibri/bCs/sz20 + green
# highlight) across all of the newly inserted text.
$rngPlaceholder = $d.Content
$foundPh = $rngPlaceholder.Find.Execute('${directorDecreto}', `
    $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($foundPh) {
    $newText = '${art8}Decreto N°140/04, del Ministerio de Salud que aprobó el Reglamento Orgánico de los Servicios de Salud, ${directorDecreto}'
    $rngPlaceholder.Text = $newText
}

# --- Step 3: the new lead-in text should NOT be highlighted -- only the
# trailing ${directorDecreto} keeps the green highlight. Locate the
# lead-in prefix precisely.
$prefixText = '${art8}Decreto N°140/04, del Ministerio de Salud que aprobó el Reglamento Orgánico de los Servicios de Salud, '
$rngPrefix = $d.Content
$foundPrefix = $rngPrefix.Find.Execute($prefixText, `
    $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($foundPrefix) {
    $ps = $rngPrefix.Start
    $pe = $rngPrefix.End

    # Assigning HighlightColorIndex on a Range produced by Find only behaves
    # (i.e. stays within the matched range) when that range is the sole
    # content of its paragraph -- otherwise it bleeds to the whole
    # paragraph. Temporarily fence the prefix off with paragraph marks,
    # clear its highlight there, then remove the fences again.
    $d.Range($pe, $pe).InsertAfter([char]13)
    $d.Range($ps, $ps).InsertBefore([char]13)

    $rngFenced = $d.Range($ps, $pe + 2)
    $rngFenced.Find.Execute($prefixText, $true, $false, $false, $false, $false, `
        $true, 1, $false, '', 0) | Out-Null
    $rngFenced.HighlightColorIndex = 0

    # Remove the two fence paragraph marks (order matters: delete the
    # earlier one first so the later offset stays valid).
    $d.Range($ps, $ps + 1).Delete()
    $d.Range($pe, $pe + 1).Delete()
}
